$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot every data cell (rows 2-41, columns A-T / 1-20) before writing anything,
# since the edit re-orders rows (a permutation) and naive in-place writes would clobber
# source rows before they are read.
$snapshot = @{}
for ($r = 2; $r -le 41; $r++) {
    for ($c = 1; $c -le 20; $c++) {
        $snapshot["$r-$c"] = $ws.Cells.Item($r, $c).Value()
    }
}

# Map: destination row -> source row (content that should end up in destination row
# comes from the source row in the ORIGINAL / before layout).
$rowMap = @{}
$rowMap[2] = 12
$rowMap[3] = 37
$rowMap[4] = 38
$rowMap[5] = 33
$rowMap[6] = 8
$rowMap[7] = 9
$rowMap[8] = 13
$rowMap[9] = 14
$rowMap[10] = 10
$rowMap[11] = 11
$rowMap[12] = 19
$rowMap[13] = 20
$rowMap[14] = 15
$rowMap[15] = 32
$rowMap[16] = 17
$rowMap[17] = 18
$rowMap[18] = 41
$rowMap[19] = 30
$rowMap[20] = 31
$rowMap[21] = 3
$rowMap[22] = 28
$rowMap[23] = 24
$rowMap[24] = 25
$rowMap[25] = 26
$rowMap[26] = 22
$rowMap[27] = 23
$rowMap[28] = 39
$rowMap[29] = 40
$rowMap[30] = 5
$rowMap[31] = 6
$rowMap[32] = 34
$rowMap[33] = 35
$rowMap[34] = 27
$rowMap[35] = 36
$rowMap[36] = 2
$rowMap[37] = 16
$rowMap[38] = 21
$rowMap[39] = 7
$rowMap[40] = 4
$rowMap[41] = 29

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $snapshot["$srcRow-$c"]
    }
}
